# Updates crypto price/volume data (and swaps the rank-25/26 rows
# between "ImmutableX" and "Monero") to match the refreshed GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.516.60"
$ws.Range("E2").Value = "  +6.60%  "
# Row 3
$ws.Range("D3").Value = "1.733.34"
$ws.Range("E3").Value = "  +4.43%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.47%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.95"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.04%  "
# Row 6
$ws.Range("E6").Value = "  +0.22%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3747"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.90%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.74"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.52%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3362"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.90%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.178"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.81%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07410"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.70%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.006"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.67%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.405"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.51%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.28"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.60%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.077"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.15%  "
# Row 16
$ws.Range("D16").Value = "1.739.03"
$ws.Range("E16").Value = "  +5.01%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001076"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.62%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06652"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.94%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.15"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.18%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.31%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.13%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.151"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.97%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.89"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.95%  "
# Row 24
$ws.Range("D24").Value = "26.551.67"
$ws.Range("E24").Value = "  +6.78%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.454"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.70%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.412"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.78%  "
# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.54"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.87%  "
# Row 28
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.385"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +15.22%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.53"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.69%  "
# Row 30
$ws.Range("D30").Value = "1.936.85"
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.66"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.15%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.155"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.55%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.015"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.27%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08599"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.72%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.700"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.12%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.84"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.63%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.385"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.13%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02335"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.18%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06244"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.85%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2164"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.56%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.528"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.15%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.228"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.60%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6183"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.21%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.16"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.03%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.002"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.29%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.905"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.54%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6001"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.67%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.63"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.98%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.046"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.00%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07228"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.56%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.27"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.17%  "
